# "sdk detailed view without EULA"
#
# Moves several boxes on the (single) slide and removes the EULA callout
# plus the two dashed "highlight" rectangles and their related "Various
# Licenses" caption.
#
# NOTE on precision: this runtime stores Shape.Left/Top/Width/Height as
# 32-bit floats (points), while the OOXML stores EMUs (1 pt = 12700 EMU).
# A naive "$emu / 12700" assignment can therefore land 1 EMU away from the
# intended value after the float32 round-trip. The literals below were
# solved so that, once truncated through float32 and multiplied back by
# 12700, they reproduce the exact target EMU offsets from the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Rounded Rectangle 114: off (6897105,2533066) -> (6897105,2348880) ---
$shp = $s.Shapes.Item("Rounded Rectangle 114")
$shp.Top = 184.95118713378906

# --- Rounded Rectangle 61: ext (2460363,1825250) -> (2460363,1472046) ---
$shp = $s.Shapes.Item("Rounded Rectangle 61")
$shp.Height = 115.90913391113281

# --- Rounded Rectangle 60: off (922655,4310180) -> (2373173,4267539) ---
$shp = $s.Shapes.Item("Rounded Rectangle 60")
$shp.Left = 186.86402893066406
$shp.Top = 336.0267028808594

# --- TextBox 63: off (7010189,1424566) -> (6984768,1077829) ---
$shp = $s.Shapes.Item("TextBox 63")
$shp.Left = 549.9817504882812
$shp.Top = 84.8684310913086

# --- TextBox 74: off (991099,4376476) -> (2441617,4333835) ---
$shp = $s.Shapes.Item("TextBox 74")
$shp.Left = 192.25331115722656
$shp.Top = 341.2468566894531

# --- TextBox 83: off (6938994,2640675) -> (6938994,2456489) ---
$shp = $s.Shapes.Item("TextBox 83")
$shp.Top = 193.42433166503906

# --- Rectangle 10 ("MICROEJ SDK End User License Agreement (EULA)") removed ---
$s.Shapes.Item("Rectangle 10").Delete()

# --- Picture 112: off (8749573,2623123) -> (8749573,2438937) ---
$shp = $s.Shapes.Item("Picture 112")
$shp.Top = 192.04229736328125

# --- TextBox 116: off (7010973,3041665) -> (7010973,2857479) ---
$shp = $s.Shapes.Item("TextBox 116")
$shp.Top = 224.99835205078125

# --- TextBox 120: off (1189031,4895808) -> (2639549,4853167) ---
$shp = $s.Shapes.Item("TextBox 120")
$shp.Left = 207.83851623535156
$shp.Top = 382.13916015625

# --- Picture 29: off (3304370,4414606) -> (4754888,4371965) ---
$shp = $s.Shapes.Item("Picture 29")
$shp.Left = 374.400634765625
$shp.Top = 344.2492370605469

# --- Rounded Rectangle 117 (dashed highlight box, no text) removed ---
$s.Shapes.Item("Rounded Rectangle 117").Delete()

# --- Rounded Rectangle 40 (dashed highlight box, no text) removed ---
$s.Shapes.Item("Rounded Rectangle 40").Delete()

# --- TextBox 33 ("Various Licenses" / "SDK EULA , Apache , Eclipse , BSD, etc.") removed ---
$s.Shapes.Item("TextBox 33").Delete()
